$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 should look just like the existing data rows (2 and 3): no
# explicit per-cell style override, just whatever the column default
# gives it. Seed row 4's formatting from row 3 before writing values so
# new cells don't pick up an inherited/explicit style stamp.
$ws.Range("A3:H3").Copy()
$ws.Range("A4:H4").PasteSpecial(-4122)

$ws.Range("A4").Value = "irBHYEJS"
$ws.Range("B4").Value = "trashboatsr"
$ws.Range("C4").Value = 1890
$ws.Range("D4").Value = 100
$ws.Range("E4").Value = "https://lichess.org/irBHYEJS"
$ws.Range("F4").Value = 4042
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = "blank"
